$wb = $excel.ActiveWorkbook

# Set the Name field on the isa_template sheet
$ws1 = $wb.Worksheets.Item("isa_template")
$ws1.Range("B3").Value = "Agronomy - Organic material applications"

# Rename Input/Output Source Name columns to Sample Name on the
# Events-OrganicFertilization sheet (this also updates the table header)
$ws2 = $wb.Worksheets.Item("Events-OrganicFertilization")
$ws2.Range("A1").Value = "Input [Sample Name]"
$ws2.Range("AD1").Value = "Output [Sample Name]"
